$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new C-column numeric values (attendance counts added next to names)
$ws.Range("C6").Value = 19
$ws.Range("C7").Value = 16
$ws.Range("C8").Value = 17
$ws.Range("C9").Value = 18
$ws.Range("C12").Value = 18
$ws.Range("C13").Value = 19
$ws.Range("C14").Value = 19
$ws.Range("C15").Value = 17
$ws.Range("C21").Value = 17
$ws.Range("C22").Value = 17
$ws.Range("C24").Value = 17

# Move the view / selection to where the author left off scrolling
$ws.Activate()
$ws.Range("E24").Select()
